$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 278.7143
$ws.Range("I41").Value = 90
$ws.Range("J41").Value = 750.5
$ws.Range("K41").Value = 90
$ws.Range("L41").Value = 750.5
$ws.Range("M41").Value = 350
$ws.Range("N41").Value = -1630.5
$ws.Range("H129").Value = 1110.3662
$ws.Range("I129").Value = 523.9
$ws.Range("J129").Value = 1206.5082
$ws.Range("K129").Value = 1571.7
$ws.Range("L129").Value = 3619.5246
$ws.Range("M129").Value = 3428.3
$ws.Range("N129").Value = -13619.5246
$ws.Range("H138").Value = 3745.3103
$ws.Range("I138").Value = 3121.348
$ws.Range("J138").Value = 3969.5469
$ws.Range("K138").Value = 9364.044
$ws.Range("L138").Value = 11908.6407
$ws.Range("M138").Value = -4224.044
$ws.Range("N138").Value = -22188.6407
$ws.Range("H140").Value = 74338.60000000001
$ws.Range("J140").Value = 74338.60000000001
$ws.Range("L140").Value = 74338.60000000001
$ws.Range("N140").Value = -84698.60000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2269
$ws.Range("I2").Value = 2252.625
$ws.Range("J2").Value = 2400
$ws.Range("K2").Value = 2252.625
$ws.Range("L2").Value = 2400
$ws.Range("M2").Value = -2139.625
$ws.Range("N2").Value = -2626
$ws.Range("H32").Value = 5272596
$ws.Range("I32").Value = 6105617.5
$ws.Range("J32").Value = 18153.846
$ws.Range("K32").Value = 6105617.5
$ws.Range("L32").Value = 18153.846
$ws.Range("M32").Value = -6105330.5
$ws.Range("N32").Value = -18727.846
$ws.Range("H97").Value = 1317.8572
$ws.Range("I97").Value = 1164.2858
$ws.Range("J97").Value = 1471.4286
$ws.Range("K97").Value = 1164.2858
$ws.Range("L97").Value = 1471.4286
$ws.Range("M97").Value = -668.2858000000001
$ws.Range("N97").Value = -2463.4286
$ws.Range("H116").Value = 2269
$ws.Range("I116").Value = 2252.625
$ws.Range("J116").Value = 2400
$ws.Range("K116").Value = 2252.625
$ws.Range("L116").Value = 2400
$ws.Range("M116").Value = 41.375
$ws.Range("N116").Value = -6988

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2269
$ws.Range("I3").Value = 2252.625
$ws.Range("J3").Value = 2400
$ws.Range("K3").Value = 2252.625
$ws.Range("L3").Value = 2400
$ws.Range("M3").Value = -2138.625
$ws.Range("N3").Value = -2628
$ws.Range("H20").Value = 35716196
$ws.Range("I20").Value = 1761.2941
$ws.Range("J20").Value = 90911224
$ws.Range("K20").Value = 1761.2941
$ws.Range("L20").Value = 90911224
$ws.Range("M20").Value = -1514.2941
$ws.Range("N20").Value = -90911718
$ws.Range("H36").Value = 20891.8
$ws.Range("I36").Value = 976.5714
$ws.Range("J36").Value = 67360.664
$ws.Range("K36").Value = 976.5714
$ws.Range("L36").Value = 67360.664
$ws.Range("M36").Value = -442.5714
$ws.Range("N36").Value = -68428.664
$ws.Range("H134").Value = 2325.6667
$ws.Range("I134").Value = 2179.432
$ws.Range("J134").Value = 3244.8572
$ws.Range("K134").Value = 6538.295999999999
$ws.Range("L134").Value = 9734.571599999999
$ws.Range("M134").Value = -4003.295999999999
$ws.Range("N134").Value = -14804.5716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 128330.95
$ws.Range("J141").Value = 124558.9
$ws.Range("L141").Value = 124558.9
$ws.Range("N141").Value = -134918.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 840.8421
$ws.Range("I68").Value = 639.5714
$ws.Range("J68").Value = 1404.4
$ws.Range("K68").Value = 1918.7142
$ws.Range("L68").Value = 4213.200000000001
$ws.Range("M68").Value = -1107.7142
$ws.Range("N68").Value = -5835.200000000001
$ws.Range("H71").Value = 840.8421
$ws.Range("I71").Value = 639.5714
$ws.Range("J71").Value = 1404.4
$ws.Range("K71").Value = 5756.1426
$ws.Range("L71").Value = 12639.6
$ws.Range("M71").Value = -1700.1426
$ws.Range("N71").Value = -20751.6
$ws.Range("H96").Value = 10000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 10000
$ws.Range("K96").Value = 0
$ws.Range("M96").Value = 30000
$ws.Range("N96").Value = -34118
$ws.Range("H113").Value = 573.75555
$ws.Range("I113").Value = 538.4286
$ws.Range("J113").Value = 697.4
$ws.Range("K113").Value = 1615.2858
$ws.Range("L113").Value = 2092.2
$ws.Range("M113").Value = 554.7142000000001
$ws.Range("N113").Value = -6432.2
$ws.Range("H131").Value = 3920.75
$ws.Range("J131").Value = 4402.2856
$ws.Range("L131").Value = 13206.8568
$ws.Range("N131").Value = -23286.8568
$ws.Range("H137").Value = 36146.605
$ws.Range("I137").Value = 6863.7144
$ws.Range("J137").Value = 87391.664
$ws.Range("K137").Value = 20591.1432
$ws.Range("L137").Value = 262174.992
$ws.Range("M137").Value = -15491.1432
$ws.Range("N137").Value = -272374.992
$ws.Range("H140").Value = 1537.2325
$ws.Range("I140").Value = 1086.0938
$ws.Range("J140").Value = 2849.6365
$ws.Range("K140").Value = 3258.2814
$ws.Range("L140").Value = 8548.9095
$ws.Range("M140").Value = 1921.7186
$ws.Range("N140").Value = -18908.9095

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1677.6
$ws.Range("I102").Value = 1643.3334
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 1643.3334
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -21.33339999999998
$ws.Range("N102").Value = -5744
$ws.Range("H122").Value = 1605.4
$ws.Range("I122").Value = 1459.8
$ws.Range("J122").Value = 1751
$ws.Range("K122").Value = 4379.4
$ws.Range("L122").Value = 5253
$ws.Range("M122").Value = -1929.4
$ws.Range("N122").Value = -10153
$ws.Range("H132").Value = 41673420
$ws.Range("I132").Value = 66675270
$ws.Range("J132").Value = 3668.3333
$ws.Range("K132").Value = 200025810
$ws.Range("L132").Value = 11004.9999
$ws.Range("M132").Value = -200023280
$ws.Range("N132").Value = -16064.9999
$ws.Range("H135").Value = 43068.57
$ws.Range("J135").Value = 43068.57
$ws.Range("L135").Value = 43068.57
$ws.Range("N135").Value = -53208.57
$ws.Range("H140").Value = 53133.332
$ws.Range("J140").Value = 53133.332
$ws.Range("L140").Value = 53133.332
$ws.Range("N140").Value = -63493.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4179.143
$ws.Range("I61").Value = 2555.889
$ws.Range("K61").Value = 2555.889
$ws.Range("M61").Value = -2353.889
$ws.Range("H100").Value = 59509.625
$ws.Range("I100").Value = 92554
$ws.Range("K100").Value = 92554
$ws.Range("M100").Value = -92013
$ws.Range("H113").Value = 4179.143
$ws.Range("I113").Value = 2555.889
$ws.Range("K113").Value = 2555.889
$ws.Range("M113").Value = -385.8890000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 457.26086
$ws.Range("I100").Value = 479.07144
$ws.Range("J100").Value = 423.33334
$ws.Range("K100").Value = 958.14288
$ws.Range("L100").Value = 846.66668
$ws.Range("M100").Value = -417.14288
$ws.Range("N100").Value = -1928.66668
$ws.Range("H107").Value = 334.33334
$ws.Range("I107").Value = 400
$ws.Range("J107").Value = 301.5
$ws.Range("K107").Value = 1200
$ws.Range("L107").Value = 904.5
$ws.Range("M107").Value = 720
$ws.Range("N107").Value = -4744.5
